# "mengubah password menjadi generatestring" - the manual "Password*" column
# is no longer needed (passwords are generated instead of entered by hand),
# so drop the whole Password column (C) from the sample-format sheet.
# Deleting the entire column shifts First Name*/Last Name*/Expired*/Quota*/Phone
# one column to the left automatically (D->C, E->D, F->E, G->F, H->G) and
# keeps the hyperlinks/styles on column B untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Columns.Item(3).Delete()

# Matches the saved selection in the target file (cell C3 instead of A3).
$ws.Range("C3").Select()
